$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (row 6) down into the new row 7
$ws.Range("A6:D6").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)

# Populate the new row's values (DIF 2020-10-30 FX entries)
$ws.Range("A7").Value = 44134
$ws.Range("B7").Value = "USD"
$ws.Range("C7").Value = "HKD"
$ws.Range("D7").Value = 7.7522

# Update the active selection to match the post-edit cursor position
$ws.Range("D8").Select()
